# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Macroferia Regional de Talca" / Piña - Caramelo
# at row 314 (pushing the previous rows 314-329 down to 315-330).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 314, shifting existing rows 314:329 down to 315:330.
$ws.Rows("314:314").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A314").Value = 5
$ws.Range("B314").Value = "Macroferia Regional de Talca"
$ws.Range("C314").Value = "Maule"
$ws.Range("D314").Value = 44939
$ws.Range("E314").Value = 7
$ws.Range("F314").Value = "Fruta"
$ws.Range("G314").Value = 100108
$ws.Range("H314").Value = "Tropicales y subtropicales"
$ws.Range("I314").Value = 100108005
$ws.Range("J314").Value = "Piña"
$ws.Range("K314").Value = "Caramelo"
$ws.Range("L314").Value = "Segunda"
$ws.Range("M314").Value = 200
$ws.Range("N314").Value = 18000
$ws.Range("O314").Value = 18000
$ws.Range("P314").Value = 18000
$ws.Range("Q314").Value = "$/caja 14 unidades"
$ws.Range("R314").Value = "Ecuador"
$ws.Range("S314").Value = 1286
$ws.Range("T314").Value = 14
